$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.999999989655598
$ws.Range("A2").Value = 0.99892520668824358
$ws.Range("A3").Value = 0.99674162412253542
$ws.Range("A4").Value = 1.0010881859886929
$ws.Range("A5").Value = 0.99234279148402416
$ws.Range("A6").Value = 0.97135695810157818
$ws.Range("A7").Value = 0.9717485591397319
$ws.Range("A8").Value = 0.97157420120646876
$ws.Range("A9").Value = 0.97603118605721129
$ws.Range("A10").Value = 0.98149058445353099
$ws.Range("A11").Value = 0.98101073881267609
$ws.Range("A12").Value = 0.98035349868552291
$ws.Range("A13").Value = 0.97849207107261338
$ws.Range("A14").Value = 0.97826723925866377
$ws.Range("A15").Value = 0.97954671633662005
$ws.Range("A16").Value = 0.98172601132029746
$ws.Range("A17").Value = 0.98742373829052299
$ws.Range("A18").Value = 0.98631487565192033
$ws.Range("A19").Value = 0.9978720682766733
$ws.Range("A20").Value = 0.99075529555898845
$ws.Range("A21").Value = 0.98935684729057538
$ws.Range("A22").Value = 0.98809235167995446
$ws.Range("A23").Value = 0.98969681736528847
$ws.Range("A24").Value = 0.98185629644646166
$ws.Range("A25").Value = 0.97539952826501919
$ws.Range("A26").Value = 0.96770484552984948
$ws.Range("A27").Value = 0.96286249921306832
$ws.Range("A28").Value = 0.9414034647519034
$ws.Range("A29").Value = 0.926137340899998
$ws.Range("A30").Value = 0.91956808960322434
$ws.Range("A31").Value = 0.911915318630517
$ws.Range("A32").Value = 0.91023608193077932
$ws.Range("A33").Value = 0.90971609097706096
